# Improve Excel cell formatting for multi-line text and HSN codes:
# - Update the HSN Codes value on row 2 (F2)
# - Append two new data rows (3 and 4) with the same look & feel as row 2
#   (row height, number/text styling, alignment) so vendor names / HSN
#   codes continue to render correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Update existing HSN Codes cell (F2)
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "996211, 300980061004, 300988526002, 300992658003, 600000562, 300922355001, 300989351001"

# ---------------------------------------------------------------------
# Helper data for the two new rows, keyed by column letter.
# Columns A (S.No.) are numeric + center/center style (same as A2).
# Column B (Vendor) and F (HSN Codes) wrap text (same as B2 / F2).
# All remaining columns are plain vertical-center text (same as C2..K2,
# excluding A/B/F) and must stay as TEXT even when the value looks like
# a number or a date, exactly like row 2's values.
# ---------------------------------------------------------------------

function Set-TextCell($cell, $srcFormatCell, $value) {
    # Clone the source cell's formatting first so the destination matches
    # row 2's look (font/fill/border/alignment) ...
    $srcFormatCell.Copy()
    $cell.PasteSpecial($xlPasteFormats)
    # ... force the value to be stored as literal text, not auto-parsed
    # into a number or date ...
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # ... then re-apply the source formatting so the NumberFormat reverts
    # back to the same style id used elsewhere in the column (General).
    $srcFormatCell.Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

function Set-NumberCell($cell, $srcFormatCell, $value) {
    $srcFormatCell.Copy()
    $cell.PasteSpecial($xlPasteFormats)
    $cell.Value = $value
}

# ---------------------------------------------------------------------
# 2. Row 3 - SONOVISION ELECTRONICS PVT LTD
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 30

Set-NumberCell $ws.Range("A3") $ws.Range("A2") 2
Set-TextCell   $ws.Range("B3") $ws.Range("B2") "SONOVISION ELECTRONICS PVT LTD"
Set-TextCell   $ws.Range("C3") $ws.Range("C2") "17/Feb/2023"
Set-TextCell   $ws.Range("D3") $ws.Range("D2") "37ABCCS7530B1ZK"
Set-TextCell   $ws.Range("E3") $ws.Range("E2") "NDYL 3826"
Set-TextCell   $ws.Range("F3") $ws.Range("F2") "85287219, 0"
Set-TextCell   $ws.Range("G3") $ws.Range("G2") "7547"
Set-TextCell   $ws.Range("H3") $ws.Range("H2") "7547"
Set-TextCell   $ws.Range("I3") $ws.Range("I2") "N/A"
Set-TextCell   $ws.Range("J3") $ws.Range("J2") "15094"
Set-TextCell   $ws.Range("K3") $ws.Range("K2") "69000"

# ---------------------------------------------------------------------
# 3. Row 4 - LAKSHMI AGENCIES
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 30

Set-NumberCell $ws.Range("A4") $ws.Range("A2") 3
Set-TextCell   $ws.Range("B4") $ws.Range("B2") "LAKSHMI AGENCIES No:18, Kannadasan Nagar Main Road, Ramapuram"
Set-TextCell   $ws.Range("C4") $ws.Range("C2") "17/07/2025"
Set-TextCell   $ws.Range("D4") $ws.Range("D2") "33AABFL7718B1ZQ"
Set-TextCell   $ws.Range("E4") $ws.Range("E2") "LA226412507098"
Set-TextCell   $ws.Range("F4") $ws.Range("F2") "15121910, 15121910, 15121910, 15180039, 15180039"
Set-TextCell   $ws.Range("G4") $ws.Range("G2") "720.00"
Set-TextCell   $ws.Range("H4") $ws.Range("H2") "720.00"
Set-TextCell   $ws.Range("I4") $ws.Range("I2") "N/A"
Set-TextCell   $ws.Range("J4") $ws.Range("J2") "1913.39"
Set-TextCell   $ws.Range("K4") $ws.Range("K2") "33725.00"
